$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated MOSIP identity schema JSON (address fields replaced with
# island/district/city/villa terms; JSON re-serialized without escaped
# forward slashes and with normalized array/object formatting).
$schemaJson = @'
{
	"$schema": "http://json-schema.org/draft-07/schema#",
	"description": "MOSIP Sample identity",
	"additionalProperties": false,
	"title": "MOSIP identity",
	"type": "object",
	"definitions": {
		"simpleType": {
			"uniqueItems": true,
			"additionalItems": false,
			"type": "array",
			"items": {
				"additionalProperties": false,
				"type": "object",
				"required": [
					"language",
					"value"
				],
				"properties": {
					"language": {
						"type": "string"
					},
					"value": {
						"type": "string"
					}
				}
			}
		},
		"documentType": {
			"additionalProperties": false,
			"type": "object",
			"properties": {
				"format": {
					"type": "string"
				},
				"type": {
					"type": "string"
				},
				"value": {
					"type": "string"
				},
				"refNumber": {
					"type": [
						"string",
						"null"
					]
				}
			}
		},
		"biometricsType": {
			"additionalProperties": false,
			"type": "object",
			"properties": {
				"format": {
					"type": "string"
				},
				"version": {
					"type": "number",
					"minimum": 0
				},
				"value": {
					"type": "string"
				}
			}
		}
	},
	"properties": {
		"identity": {
			"additionalProperties": false,
			"type": "object",
			"required": [
				"IDSchemaVersion",
				"fullName",
				"dateOfBirth",
				"gender",
				"island",
				"district",
				"city",
				"villa",
				"individualBiometrics"
			],
			"properties": {
				"proofOfAddress": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/documentType"
				},
				"gender": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"city": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^(?=.{0,50}$).*",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"proofOfException-1": {
					"bioAttributes": [],
					"fieldCategory": "evidence",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/documentType"
				},
				"referenceIdentityNumber": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^([0-9]{10,30})$",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "kyc",
					"type": "string",
					"fieldType": "default"
				},
				"individualBiometrics": {
					"bioAttributes": [
						"leftEye",
						"rightEye",
						"rightIndex",
						"rightLittle",
						"rightRing",
						"rightMiddle",
						"leftIndex",
						"leftLittle",
						"leftRing",
						"leftMiddle",
						"leftThumb",
						"rightThumb",
						"face"
					],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/biometricsType"
				},
				"island": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^(?=.{0,50}$).*",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"district": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"city": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"villa": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"proofOfDateOfBirth": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/documentType"
				},
				"addressLine1": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^(?=.{0,50}$).*",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"email": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^[A-Za-z0-9_\\-]+(\\.[A-Za-z0-9_]+)*@[A-Za-z0-9_-]+(\\.[A-Za-z0-9_]+)*(\\.[a-zA-Z]{2,})$",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "none",
					"type": "string",
					"fieldType": "default"
				},
				"introducerRID": {
					"bioAttributes": [],
					"fieldCategory": "evidence",
					"format": "none",
					"type": "string",
					"fieldType": "default"
				},
				"introducerBiometrics": {
					"bioAttributes": [
						"leftEye",
						"rightEye",
						"rightIndex",
						"rightLittle",
						"rightRing",
						"rightMiddle",
						"leftIndex",
						"leftLittle",
						"leftRing",
						"leftMiddle",
						"leftThumb",
						"rightThumb",
						"face"
					],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/biometricsType"
				},
				"fullName": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^(?=.{3,50}$).*",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"dateOfBirth": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^(1869|18[7-9][0-9]|19[0-9][0-9]|20[0-9][0-9])/([0][1-9]|1[0-2])/([0][1-9]|[1-2][0-9]|3[01])$",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "none",
					"type": "string",
					"fieldType": "default"
				},
				"individualAuthBiometrics": {
					"bioAttributes": [
						"leftEye",
						"rightEye",
						"rightIndex",
						"rightLittle",
						"rightRing",
						"rightMiddle",
						"leftIndex",
						"leftLittle",
						"leftRing",
						"leftMiddle",
						"leftThumb",
						"rightThumb",
						"face"
					],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/biometricsType"
				},
				"introducerUIN": {
					"bioAttributes": [],
					"fieldCategory": "evidence",
					"format": "none",
					"type": "string",
					"fieldType": "default"
				},
				"proofOfIdentity": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/documentType"
				},
				"IDSchemaVersion": {
					"bioAttributes": [],
					"fieldCategory": "none",
					"format": "none",
					"type": "number",
					"fieldType": "default",
					"minimum": 0
				},
				"proofOfException": {
					"bioAttributes": [],
					"fieldCategory": "evidence",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/documentType"
				},
				"phone": {
					"bioAttributes": [],
					"validators": [
						{
							"validator": "^[+]*([0-9]{1})([0-9]{9})$",
							"arguments": [],
							"type": "regex"
						}
					],
					"fieldCategory": "pvt",
					"format": "none",
					"type": "string",
					"fieldType": "default"
				},
				"introducerName": {
					"bioAttributes": [],
					"fieldCategory": "evidence",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/simpleType"
				},
				"proofOfRelationship": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"fieldType": "default",
					"$ref": "#/definitions/documentType"
				},
				"UIN": {
					"bioAttributes": [],
					"fieldCategory": "none",
					"format": "none",
					"type": "string",
					"fieldType": "default"
				},
				"preferredLang": {
					"bioAttributes": [],
					"fieldCategory": "pvt",
					"format": "none",
					"type": "string",
					"fieldType": "dynamic"
				}
			}
		}
	}
}
'@

$ws.Range("F2").Value = $schemaJson

# The row holding the schema JSON keeps its "auto" max height in the
# newer authoring tool's output (409.5 -> 409.6).
$ws.Rows.Item(2).RowHeight = 409.6
